$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text entry (leading apostrophe) so numeric-looking values
# like "1.032" or "27.613.66" stay strings, matching the source data.

$ws.Range("D2").Value = "'27.613.66"
$ws.Range("E2").Value = "'  +2.94%  "
$ws.Range("D3").Value = "'1.851.00"
$ws.Range("E3").Value = "'  +2.08%  "
$ws.Range("D4").Value = "'1.032"
$ws.Range("E4").Value = "'  +2.70%  "
$ws.Range("D5").Value = "'321.46"
$ws.Range("E5").Value = "'  +4.00%  "
$ws.Range("D6").Value = "'1.028"
$ws.Range("E6").Value = "'  +2.38%  "
$ws.Range("D7").Value = "'0.4374"
$ws.Range("E7").Value = "'  +1.13%  "
$ws.Range("D8").Value = "'0.3752"
$ws.Range("E8").Value = "'  +1.15%  "
$ws.Range("D9").Value = "'0.07400"
$ws.Range("E9").Value = "'  +1.98%  "
$ws.Range("D10").Value = "'0.8752"
$ws.Range("E10").Value = "'  +0.97%  "
$ws.Range("E11").Value = "'  +2.80%  "
$ws.Range("D12").Value = "'1.855.72"
$ws.Range("E12").Value = "'  -4.80%  "
$ws.Range("D13").Value = "'5.513"
$ws.Range("E13").Value = "'  +2.88%  "
$ws.Range("D14").Value = "'6.690"
$ws.Range("E14").Value = "'  +0.70%  "
$ws.Range("D15").Value = "'0.07180"
$ws.Range("D16").Value = "'82.66"
$ws.Range("E16").Value = "'  +2.54%  "
$ws.Range("D17").Value = "'1.033"
$ws.Range("E17").Value = "'  +2.80%  "
$ws.Range("D18").Value = "'0.000009042"
$ws.Range("E18").Value = "'  +1.88%  "
$ws.Range("D19").Value = "'1.027"
$ws.Range("E19").Value = "'  +2.24%  "
$ws.Range("D20").Value = "'15.43"
$ws.Range("E20").Value = "'  +1.37%  "
$ws.Range("D21").Value = "'27.613.90"
$ws.Range("E21").Value = "'  +2.82%  "
$ws.Range("D22").Value = "'5.259"
$ws.Range("E22").Value = "'  +0.81%  "
$ws.Range("D23").Value = "'11.23"
$ws.Range("E23").Value = "'  +0.34%  "
$ws.Range("D24").Value = "'2.075.00"
$ws.Range("E24").Value = "'  -4.57%  "
$ws.Range("D25").Value = "'157.51"
$ws.Range("E25").Value = "'  +2.25%  "
$ws.Range("D26").Value = "'1.942"
$ws.Range("E26").Value = "'  +3.81%  "
$ws.Range("D27").Value = "'18.74"
$ws.Range("E27").Value = "'  +2.46%  "
$ws.Range("D28").Value = "'5.299"
$ws.Range("E28").Value = "'  +1.31%  "
$ws.Range("D29").Value = "'1.935"
$ws.Range("E29").Value = "'  +1.80%  "
$ws.Range("D30").Value = "'116.13"
$ws.Range("E30").Value = "'  +0.66%  "
$ws.Range("D31").Value = "'0.09065"
$ws.Range("E31").Value = "'  +1.34%  "
$ws.Range("D32").Value = "'1.208"
$ws.Range("E32").Value = "'  +2.89%  "
$ws.Range("D33").Value = "'0.7677"
$ws.Range("E33").Value = "'  +1.48%  "
$ws.Range("E34").Value = "'  +1.96%  "
$ws.Range("D35").Value = "'2.877"
$ws.Range("E35").Value = "'  +2.38%  "
$ws.Range("E36").Value = "'  +1.98%  "
$ws.Range("D37").Value = "'1.154"
$ws.Range("E37").Value = "'  +1.91%  "
$ws.Range("D38").Value = "'0.01978"
$ws.Range("E38").Value = "'  +2.61%  "
$ws.Range("D39").Value = "'0.05284"
$ws.Range("E39").Value = "'  +0.94%  "
$ws.Range("D40").Value = "'2.822"
$ws.Range("E40").Value = "'  +4.82%  "
$ws.Range("D41").Value = "'0.5178"
$ws.Range("E41").Value = "'  +1.79%  "
$ws.Range("D42").Value = "'0.1673"
$ws.Range("E42").Value = "'  +1.42%  "
$ws.Range("D43").Value = "'6.736"
$ws.Range("E43").Value = "'  +2.77%  "
$ws.Range("D44").Value = "'8.599"
$ws.Range("B45").Value = "'EnergySwap"
$ws.Range("C45").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'10.65"
$ws.Range("E45").Value = "'  +1.70%  "
$ws.Range("B46").Value = "'Quant"
$ws.Range("C46").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "'108.86"
$ws.Range("E46").Value = "'  +1.92%  "
$ws.Range("D47").Value = "'1.713"
$ws.Range("E47").Value = "'  +3.39%  "
$ws.Range("D48").Value = "'0.4657"
$ws.Range("E48").Value = "'  +1.60%  "
$ws.Range("D49").Value = "'0.06389"
$ws.Range("D50").Value = "'1.879"
$ws.Range("E50").Value = "'  +3.59%  "
$ws.Range("E51").Value = "'  +5.20%  "
